$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.104631
$ws.Range("H2").Value = 3.313893
$ws.Range("I2").Value = 0.8734451962653081
$ws.Range("J2").Value = 0.8734451962653083
$ws.Range("M2").Value = 1.704784666666667
$ws.Range("N2").Value = 5.114354000000001
$ws.Range("O2").Value = 0.981710618882129
$ws.Range("P2").Value = 0.981710618882129
$ws.Range("Q2").Value = 1.883157991124667
$ws.Range("R2").Value = 16.948421920122
$ws.Range("S2").Value = 0.8574704241852382
$ws.Range("T2").Value = 0.8574704241852383

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.104631
$ws.Range("H3").Value = 3.313893
$ws.Range("I3").Value = 0.8734451962653081
$ws.Range("J3").Value = 0.8734451962653083
$ws.Range("N3").Value = 0.095281
$ws.Range("O3").Value = 0.01828938111787102
$ws.Range("P3").Value = 0.01828938111787102
$ws.Range("Q3").Value = 0.03508344877033334
$ws.Range("R3").Value = 0.315751038933
$ws.Range("S3").Value = 0.01597477208006987
$ws.Range("T3").Value = 0.01597477208006988

# Row 4
$ws.Range("I4").Value = 0.1265548037346918
$ws.Range("J4").Value = 0.1265548037346918
$ws.Range("M4").Value = 1.704784666666667
$ws.Range("N4").Value = 5.114354000000001
$ws.Range("O4").Value = 0.981710618882129
$ws.Range("P4").Value = 0.981710618882129
$ws.Range("Q4").Value = 0.2728536272077778
$ws.Range("R4").Value = 2.45568264487
$ws.Range("S4").Value = 0.1242401946968907
$ws.Range("T4").Value = 0.1242401946968907

# Row 5
$ws.Range("I5").Value = 0.1265548037346918
$ws.Range("J5").Value = 0.1265548037346918
$ws.Range("N5").Value = 0.095281
$ws.Range("O5").Value = 0.01828938111787102
$ws.Range("P5").Value = 0.01828938111787102
$ws.Range("Q5").Value = 0.005083294283888889
$ws.Range("R5").Value = 0.045749648555
$ws.Range("S5").Value = 0.002314609037801145
$ws.Range("T5").Value = 0.002314609037801145
